# Update the ranking worksheet: matrices scores were refreshed and the
# table re-sorted (descending by matrices score), causing some rows to
# swap places while each person's own identity (index/prolificid/name/
# gender) travels with them. mat_rank (column H) and the raw row index
# (column A) are unaffected by the resort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, index, prolificid, name, gender, matrices
$rows = @(
    @{ Row = 2;  Idx = 0;  Prolific = "5eeaa065c7acf61c4322f6d9"; Name = "Yonifredy"; Gender = "male";   Matrices = 15.02328293437414 },
    @{ Row = 3;  Idx = 1;  Prolific = "5e0adc8f4cac6834756db412"; Name = "Mary";      Gender = "female"; Matrices = 11.12005548300506 },
    @{ Row = 4;  Idx = 3;  Prolific = "60ba8ba51a5e0a105396888a"; Name = "Alfredo";   Gender = "male";   Matrices = 10.3560449567461 },
    @{ Row = 5;  Idx = 2;  Prolific = "60778ed0fde3e9c3a96f1d11"; Name = "Melissa";   Gender = "female"; Matrices = 10.15590669353794 },
    @{ Row = 6;  Idx = 4;  Prolific = "5e706891c396cc64388ef760"; Name = "Maria";     Gender = "female"; Matrices = 9.075645813370125 },
    @{ Row = 7;  Idx = 6;  Prolific = "5dd671942b033b5ec8bc97b4"; Name = "Juan";      Gender = "male";   Matrices = 7.229575176107406 },
    @{ Row = 8;  Idx = 7;  Prolific = "6024c18b094ac71dd93f4f5a"; Name = "Katherine"; Gender = "female"; Matrices = 5.144726965691964 },
    @{ Row = 9;  Idx = 9;  Prolific = "5e35d91ea42bce592e996843"; Name = "Sergio";    Gender = "male";   Matrices = 5.106254872490608 },
    @{ Row = 10; Idx = 8;  Prolific = "5f0142aa1eb1e528e7abce50"; Name = "Valeria";   Gender = "female"; Matrices = 5.051234491524045 },
    @{ Row = 11; Idx = 10; Prolific = "60743a8fd12c5ffa72972fd5"; Name = "Josue";     Gender = "male";   Matrices = 4.078136080597864 },
    @{ Row = 12; Idx = 12; Prolific = "5e58b3e415b8d40b5e1dabf1"; Name = "Cristian";  Gender = "male";   Matrices = 3.427904729701768 },
    @{ Row = 13; Idx = 11; Prolific = "5f5ea8227fa75676f56f9276"; Name = "Carlos";    Gender = "male";   Matrices = 3.301880844181574 }
)

foreach ($r in $rows) {
    $ws.Range("B$($r.Row)").Value = $r.Idx
    $ws.Range("C$($r.Row)").Value = $r.Prolific
    $ws.Range("D$($r.Row)").Value = $r.Name
    $ws.Range("E$($r.Row)").Value = $r.Gender
    $ws.Range("F$($r.Row)").Value = $r.Matrices
}
